$wb = $excel.ActiveWorkbook

# Add the new "OOTB Domain Groups" worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "OOTB Domain Groups"

$data = New-Object 'object[,]' 14,142
$data[0,0] = "Gmail"
$data[0,1] = "gmail.com"
$data[0,2] = "googlemail.com"
$data[0,3] = "googlemail.co.uk"
$data[1,0] = "Microsoft"
$data[1,1] = "live.com"
$data[1,2] = "msn.com"
$data[1,3] = "hotmail.ca"
$data[1,4] = "hotmail.com"
$data[1,5] = "hotmail.de"
$data[1,6] = "hotmail.dk"
$data[1,7] = "hotmail.co.jp"
$data[1,8] = "hotmail.it"
$data[1,9] = "hotmail.es"
$data[1,10] = "hotmail.fr"
$data[1,11] = "hotmail.co.uk"
$data[1,12] = "hotmail.co.kr"
$data[1,13] = "hotmail.com.au"
$data[1,14] = "hotmail.com.ar"
$data[1,15] = "hotmail.co.il"
$data[1,16] = "hotmail.com.br"
$data[1,17] = "hotmail.com.tr"
$data[1,18] = "hotmail.co.th"
$data[1,19] = "hotmail.jp"
$data[1,20] = "hotmail.se"
$data[1,21] = "live.at"
$data[1,22] = "live.be"
$data[1,23] = "live.ca"
$data[1,24] = "live.cl"
$data[1,25] = "live.cn"
$data[1,26] = "live.co.kr"
$data[1,27] = "live.co.uk"
$data[1,28] = "live.com.ar"
$data[1,29] = "live.com.au"
$data[1,30] = "live.com.mx"
$data[1,31] = "live.com.my"
$data[1,32] = "live.com.sg"
$data[1,33] = "live.de"
$data[1,34] = "live.dk"
$data[1,35] = "live.fr"
$data[1,36] = "live.hk"
$data[1,37] = "live.ie"
$data[1,38] = "live.in"
$data[1,39] = "live.it"
$data[1,40] = "live.jp"
$data[1,41] = "live.nl"
$data[1,42] = "live.no"
$data[1,43] = "live.ru"
$data[1,44] = "live.se"
$data[1,45] = "outlook.com"
$data[1,46] = "live.co.uk"
$data[1,47] = "hotmail.gr"
$data[1,48] = "windowslive.com"
$data[1,49] = "xbox.com"
$data[1,50] = "hotmail.cl"
$data[1,51] = "live.at"
$data[1,52] = "live.jp"
$data[1,53] = "live.ca"
$data[1,54] = "hotmail.ca"
$data[1,55] = "hotmail.ch"
$data[1,56] = "live.fr"
$data[1,57] = "live.it"
$data[1,58] = "live.nl"
$data[1,59] = "outlook.ie"
$data[1,60] = "outlook.com.br"
$data[1,61] = "live.com.pt"
$data[1,62] = "live.be"
$data[1,63] = "live.co.za"
$data[1,64] = "mts.net"
$data[2,0] = "Yahoo"
$data[2,1] = "yahoo.com"
$data[2,2] = "rocketmail.com"
$data[2,3] = "rogers.com"
$data[2,4] = "sky.com"
$data[2,5] = "talk21.com"
$data[2,6] = "y7mail.com"
$data[2,7] = "yahoo.at"
$data[2,8] = "yahoo.be"
$data[2,9] = "yahoo.bg"
$data[2,10] = "yahoo.ca"
$data[2,11] = "yahoo.cl"
$data[2,12] = "yahoo.co.hu"
$data[2,13] = "yahoo.co.id"
$data[2,14] = "yahoo.co.il"
$data[2,15] = "yahoo.co.in"
$data[2,16] = "yahoo.co.jp"
$data[2,17] = "yahoo.co.kr"
$data[2,18] = "yahoo.com.ar"
$data[2,19] = "yahoo.com.au"
$data[2,20] = "yahoo.com.biz"
$data[2,21] = "yahoo.com.br"
$data[2,22] = "yahoo.com.cn"
$data[2,23] = "yahoo.com.co"
$data[2,24] = "yahoo.com.hk"
$data[2,25] = "yahoo.com.hr"
$data[2,26] = "yahoo.com.in"
$data[2,27] = "yahoo.com.jp"
$data[2,28] = "yahoo.com.kr"
$data[2,29] = "yahoo.com.mx"
$data[2,30] = "yahoo.com.my"
$data[2,31] = "yahoo.com.net"
$data[2,32] = "yahoo.com.pe"
$data[2,33] = "yahoo.com.ph"
$data[2,34] = "yahoo.com.sg"
$data[2,35] = "yahoo.com.tr"
$data[2,36] = "yahoo.com.tw"
$data[2,37] = "yahoo.com.ua"
$data[2,38] = "yahoo.com.ve"
$data[2,39] = "yahoo.com.vn"
$data[2,40] = "yahoo.co.nz"
$data[2,41] = "yahoo.co.th"
$data[2,42] = "yahoo.co.uk"
$data[2,43] = "yahoo.co.za"
$data[2,44] = "yahoo.cz"
$data[2,45] = "yahoo.de"
$data[2,46] = "yahoo.dk"
$data[2,47] = "yahoo.ee"
$data[2,48] = "yahoo.es"
$data[2,49] = "yahoo.fi"
$data[2,50] = "yahoo.fr"
$data[2,51] = "yahoogroups.co.kr"
$data[2,52] = "yahoogroups.com.cn"
$data[2,53] = "yahoogroups.com.sg"
$data[2,54] = "yahoogroups.com.tw"
$data[2,55] = "yahoogrupper.dk"
$data[2,56] = "yahoogruppi.it"
$data[2,57] = "yahoo.gr"
$data[2,58] = "yahoo.hr"
$data[2,59] = "yahoo.hu"
$data[2,60] = "yahoo.ie"
$data[2,61] = "yahoo.in"
$data[2,62] = "yahoo.it"
$data[2,63] = "yahoo.lt"
$data[2,64] = "yahoo.lv"
$data[2,65] = "yahoo.nl"
$data[2,66] = "yahoo.no"
$data[2,67] = "yahoo.pl"
$data[2,68] = "yahoo.pt"
$data[2,69] = "yahoo.ro"
$data[2,70] = "yahoo.rs"
$data[2,71] = "yahoo.se"
$data[2,72] = "yahoo.si"
$data[2,73] = "yahoo.sk"
$data[2,74] = "yahooxtra.co.nz"
$data[2,75] = "ymail.com"
$data[2,76] = "aol.com"
$data[2,77] = "aim.com"
$data[2,78] = "compuserve.com"
$data[2,79] = "cs.com"
$data[2,80] = "netscape.com"
$data[2,81] = "netscape.net"
$data[2,82] = "wmconnect.com"
$data[2,83] = "aol.co.uk"
$data[2,84] = "aol.in"
$data[2,85] = "aol.de"
$data[2,86] = "aol.fr"
$data[2,87] = "aol.nl"
$data[2,88] = "aol.pl"
$data[2,89] = "aol.jp"
$data[2,90] = "aol.es"
$data[2,91] = "aol.it"
$data[2,92] = "aol.com.ar"
$data[2,93] = "aol.fi"
$data[2,94] = "aol.cl"
$data[2,95] = "aol.com.co"
$data[2,96] = "aol.com.ve"
$data[2,97] = "aol.com.au"
$data[2,98] = "aol.at"
$data[2,99] = "aol.be"
$data[2,100] = "aol.com.br"
$data[2,101] = "aol.cz"
$data[2,102] = "aol.dk"
$data[2,103] = "myaol.jp"
$data[2,104] = "aolnorge.no"
$data[2,105] = "aolpolska.pl"
$data[2,106] = "aolpolcka.pl"
$data[2,107] = "aolpoland.pl"
$data[2,108] = "aol.ru"
$data[2,109] = "aol.kr"
$data[2,110] = "aol.se"
$data[2,111] = "aol.ch"
$data[2,112] = "aol.com.tr"
$data[2,113] = "aol.co.nz"
$data[2,114] = "aolchina.com"
$data[2,115] = "aol.hk"
$data[2,116] = "aol.tw"
$data[2,117] = "luckymail.com"
$data[2,118] = "verizon.net"
$data[2,119] = "aol.com.mx"
$data[2,120] = "bellatlantic.net"
$data[2,121] = "citlink.net"
$data[2,122] = "frontier.com"
$data[2,123] = "frontiernet.net"
$data[2,124] = "games.com"
$data[2,125] = "goowy.com"
$data[2,126] = "gte.net"
$data[2,127] = "love.com"
$data[2,128] = "verizon.net.in"
$data[2,129] = "wild4music.com"
$data[2,130] = "wow.com"
$data[2,131] = "yahoo.cn"
$data[2,132] = "yahoo.ne.jp"
$data[2,133] = "yahoogroups.ca"
$data[2,134] = "yahoogroups.co.in"
$data[2,135] = "yahoogroups.co.uk"
$data[2,136] = "yahoogroups.com"
$data[2,137] = "yahoogroups.com.au"
$data[2,138] = "yahoogroups.com.hk"
$data[2,139] = "yahoogroups.de"
$data[2,140] = "ybb.ne.jp"
$data[2,141] = "ygm.com"
$data[3,0] = "Apple"
$data[3,1] = "mac.com"
$data[3,2] = "icloud.com"
$data[3,3] = "apple.com"
$data[3,4] = "me.com"
$data[4,0] = "Comcast"
$data[4,1] = "comcast.net"
$data[5,0] = "Orange"
$data[5,1] = "orange.fr"
$data[5,2] = "orange.com"
$data[5,3] = "wanadoo.fr"
$data[5,4] = "francetelecom.com"
$data[5,5] = "voila.fr"
$data[5,6] = "voila.com"
$data[6,0] = "La Poste"
$data[6,1] = "laposte.net"
$data[7,0] = "Italia Online"
$data[7,1] = "libero.it"
$data[7,2] = "inwind.it"
$data[7,3] = "iol.it"
$data[7,4] = "blu.it"
$data[7,5] = "giallo.it"
$data[7,6] = "virgilio.it"
$data[8,0] = "WP"
$data[8,1] = "wp.pl"
$data[8,2] = "o2.pl"
$data[9,0] = "United Internet"
$data[9,1] = "web.de"
$data[9,2] = "gmx.de"
$data[9,3] = "gmx.ch"
$data[9,4] = "gmx.net"
$data[9,5] = "gmx.com"
$data[9,6] = "gmx.at"
$data[9,7] = "gmx.fr"
$data[9,8] = "mail.com"
$data[9,9] = "1and1.com"
$data[9,10] = "1und1.de"
$data[10,0] = "Bigpond"
$data[10,1] = "bigpond.com"
$data[10,2] = "bigpond.net.au"
$data[10,3] = "bigpond.com.au"
$data[10,4] = "telstra.com"
$data[10,5] = "bigpond.net"
$data[11,0] = "Docomo"
$data[11,1] = "docomo.ne.jp"
$data[12,0] = "Softbank"
$data[12,1] = "softbank.ne.jp"
$data[12,2] = "c.vodafone.ne.jp"
$data[12,3] = "d.vodafone.ne.jp"
$data[12,4] = "h.vodafone.ne.jp"
$data[12,5] = "k.vodafone.ne.jp"
$data[12,6] = "n.vodafone.ne.jp"
$data[12,7] = "q.vodafone.ne.jp"
$data[12,8] = "r.vodafone.ne.jp"
$data[12,9] = "s.vodafone.ne.jp"
$data[12,10] = "t.vodafone.ne.jp"
$data[12,11] = "jp-c.ne.jp"
$data[12,12] = "jp-d.ne.jp"
$data[12,13] = "jp-h.ne.jp"
$data[12,14] = "jp-k.ne.jp"
$data[12,15] = "jp-n.ne.jp"
$data[12,16] = "jp-q.ne.jp"
$data[12,17] = "jp-r.ne.jp"
$data[12,18] = "jp-s.ne.jp"
$data[12,19] = "jp-t.ne.jp"
$data[13,0] = "KDDI"
$data[13,1] = "au.com"
$data[13,2] = "ezweb.ne.jp"
$data[13,3] = "uqmobile.jp"

$ws3.Range("A1:EL14").Value = $data

# Selection bookkeeping for the new sheet (matches author's saved view state).
$ws3.Range("A23").Select()

# Restore "Warmup Plan" as the active sheet/tab, matching the source workbook.
$wb.Worksheets.Item("Warmup Plan").Activate()
$wb.Worksheets.Item("Warmup Plan").Range("A1").Select()
